# Highlight specific SLL-assignment bullet items in red, per the
# "13 September 2024/ SLL Questions" commit.
#
# wdRed = 6 (WdColorIndex)
$wdRed = 6

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    # Paragraph.Range.Text includes the trailing paragraph-mark character
    # (CR, chr 13); strip it so exact-text comparisons work.
    $t = $p.Range.Text.TrimEnd([char]13)

    if ($t -eq "Remove duplicate nodes from a sorted linked list.") {
        # Only the run text is highlighted; the paragraph mark itself is
        # left untouched.
        $p.Range.HighlightColorIndex = $wdRed
    }
    elseif ($t -eq "Remove duplicate nodes from an unsorted linked list.") {
        # Both the paragraph mark and the run get highlighted.
        $p.Range.Font.HighlightColorIndex = $wdRed
    }
    elseif ($t -eq "Splitting the List: Split the linked list into two halves.") {
        # Paragraph mark plus every run in the paragraph.
        $p.Range.Font.HighlightColorIndex = $wdRed
    }
    elseif ($t -eq "Deleting the Entire List: Delete all nodes and free up memory.") {
        # Paragraph mark plus every run in the paragraph.
        $p.Range.Font.HighlightColorIndex = $wdRed
    }
}
